# Auto-generated edit script: update cached market-price / profit figures
# across the ALC, ARM, CRP, CUL, LTW and WVR sheets (scheduled-runner refresh).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3600
$ws.Range("I51").Value = 2166.6667
$ws.Range("J51").Value = 4316.6665
$ws.Range("K51").Value = 2166.6667
$ws.Range("L51").Value = 4316.6665
$ws.Range("M51").Value = -1682.6667
$ws.Range("N51").Value = -5284.6665
$ws.Range("H99").Value = 169289.83
$ws.Range("I99").Value = 186
$ws.Range("J99").Value = 338393.66
$ws.Range("K99").Value = 558
$ws.Range("L99").Value = 1015180.98
$ws.Range("M99").Value = 940
$ws.Range("N99").Value = -1018176.98
$ws.Range("H113").Value = 2560.7827
$ws.Range("I113").Value = 2571.4285
$ws.Range("K113").Value = 2571.4285
$ws.Range("M113").Value = 682.5715
$ws.Range("H132").Value = 2231.2285
$ws.Range("I132").Value = 2168.5254
$ws.Range("J132").Value = 2567.5454
$ws.Range("K132").Value = 6505.5762
$ws.Range("L132").Value = 7702.6362
$ws.Range("M132").Value = -3975.5762
$ws.Range("N132").Value = -12762.6362
$ws.Range("H135").Value = 2249.375
$ws.Range("I135").Value = 1222.6666
$ws.Range("J135").Value = 5329.5
$ws.Range("K135").Value = 11003.9994
$ws.Range("L135").Value = 47965.5
$ws.Range("M135").Value = -8468.999400000001
$ws.Range("N135").Value = -53035.5
$ws.Range("H138").Value = 1703.8223
$ws.Range("I138").Value = 1459.697
$ws.Range("J138").Value = 2375.1667
$ws.Range("K138").Value = 4379.090999999999
$ws.Range("L138").Value = 7125.500100000001
$ws.Range("M138").Value = 760.9090000000006
$ws.Range("N138").Value = -17405.5001
$ws.Range("H139").Value = 49925
$ws.Range("J139").Value = 49925
$ws.Range("L139").Value = 49925
$ws.Range("N139").Value = -60205
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7354064.5
$ws.Range("I61").Value = 8772998
$ws.Range("J61").Value = 1409.3636
$ws.Range("K61").Value = 8772998
$ws.Range("L61").Value = 1409.3636
$ws.Range("M61").Value = -8772786
$ws.Range("N61").Value = -1833.3636
$ws.Range("H63").Value = 3357.8572
$ws.Range("I63").Value = 3251.25
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 3251.25
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -2565.25
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 3357.8572
$ws.Range("I66").Value = 3251.25
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 16256.25
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -12824.25
$ws.Range("N66").Value = -24364
$ws.Range("H74").Value = 10418131
$ws.Range("I74").Value = 13159251
$ws.Range("J74").Value = 1872.8
$ws.Range("K74").Value = 13159251
$ws.Range("L74").Value = 1872.8
$ws.Range("M74").Value = -13158377
$ws.Range("N74").Value = -3620.8
$ws.Range("H77").Value = 10418131
$ws.Range("I77").Value = 13159251
$ws.Range("J77").Value = 1872.8
$ws.Range("K77").Value = 65796255
$ws.Range("L77").Value = 9364
$ws.Range("M77").Value = -65791887
$ws.Range("N77").Value = -18100
$ws.Range("H97").Value = 11828.091
$ws.Range("I97").Value = 12988.625
$ws.Range("K97").Value = 12988.625
$ws.Range("M97").Value = -12492.625
$ws.Range("H136").Value = 7354064.5
$ws.Range("I136").Value = 8772998
$ws.Range("J136").Value = 1409.3636
$ws.Range("K136").Value = 26318994
$ws.Range("L136").Value = 4228.0908
$ws.Range("M136").Value = -26316444
$ws.Range("N136").Value = -9328.0908

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 608.9091
$ws.Range("I22").Value = 320
$ws.Range("J22").Value = 849.6667
$ws.Range("K22").Value = 320
$ws.Range("L22").Value = 849.6667
$ws.Range("M22").Value = 30
$ws.Range("N22").Value = -1549.6667
$ws.Range("H31").Value = 6062655
$ws.Range("I31").Value = 1705.325
$ws.Range("J31").Value = 22225186
$ws.Range("K31").Value = 1705.325
$ws.Range("L31").Value = 22225186
$ws.Range("M31").Value = -1410.325
$ws.Range("N31").Value = -22225776
$ws.Range("H34").Value = 6062655
$ws.Range("I34").Value = 1705.325
$ws.Range("J34").Value = 22225186
$ws.Range("K34").Value = 1705.325
$ws.Range("L34").Value = 22225186
$ws.Range("M34").Value = -1503.325
$ws.Range("N34").Value = -22225590
$ws.Range("H38").Value = 17369
$ws.Range("I38").Value = 2500
$ws.Range("J38").Value = 22325.334
$ws.Range("K38").Value = 2500
$ws.Range("L38").Value = 22325.334
$ws.Range("M38").Value = -2123
$ws.Range("N38").Value = -23079.334
$ws.Range("H46").Value = 17369
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 22325.334
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 22325.334
$ws.Range("M46").Value = -2289
$ws.Range("N46").Value = -22747.334
$ws.Range("H81").Value = 48926.285
$ws.Range("J81").Value = 48926.285
$ws.Range("L81").Value = 48926.285
$ws.Range("N81").Value = -50922.285
$ws.Range("H84").Value = 48926.285
$ws.Range("J84").Value = 48926.285
$ws.Range("L84").Value = 146778.855
$ws.Range("N84").Value = -156762.855
$ws.Range("H132").Value = 8334958.5
$ws.Range("I132").Value = 11629287
$ws.Range("J132").Value = 2244.4707
$ws.Range("K132").Value = 34887861
$ws.Range("L132").Value = 6733.4121
$ws.Range("M132").Value = -34885331
$ws.Range("N132").Value = -11793.4121
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 731073.8
$ws.Range("I2").Value = 116.6
$ws.Range("J2").Value = 1543248.5
$ws.Range("K2").Value = 699.5999999999999
$ws.Range("L2").Value = 9259491
$ws.Range("M2").Value = -586.5999999999999
$ws.Range("N2").Value = -9259717
$ws.Range("H125").Value = 3749.125
$ws.Range("J125").Value = 4113.2856
$ws.Range("L125").Value = 12339.8568
$ws.Range("N125").Value = -22179.8568

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7600
$ws.Range("I40").Value = 25500
$ws.Range("J40").Value = 3125
$ws.Range("K40").Value = 25500
$ws.Range("L40").Value = 3125
$ws.Range("M40").Value = -25364
$ws.Range("N40").Value = -3397
$ws.Range("H136").Value = 10419856
$ws.Range("I136").Value = 12500939
$ws.Range("J136").Value = 14438.125
$ws.Range("K136").Value = 37502817
$ws.Range("L136").Value = 43314.375
$ws.Range("M136").Value = -37500267
$ws.Range("N136").Value = -48414.375
$ws.Range("H139").Value = 55625
$ws.Range("J139").Value = 55625
$ws.Range("L139").Value = 55625
$ws.Range("N139").Value = -65905

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 933.1539
$ws.Range("I113").Value = 392.33334
$ws.Range("J113").Value = 2150
$ws.Range("K113").Value = 1177.00002
$ws.Range("L113").Value = 6450
$ws.Range("M113").Value = 992.9999800000001
$ws.Range("N113").Value = -10790
$ws.Range("H138").Value = 45964.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 45964.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 45964.5
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -56244.5

Write-Host "Applied scheduled-runner Sheets update."